# Atualização automática de pedidos - 30/05/2025 08:59
# Fixes D10 (RACK/seccao) on both sheets to a true numeric value, then
# appends two new order rows (REQ-010, REQ-011) to "Pedidos" and "Itens".

$wb = $excel.ActiveWorkbook
$wsPedidos = $wb.Worksheets.Item("Pedidos")
$wsItens   = $wb.Worksheets.Item("Itens")

# --- Pedidos: correct D10 ("RACK") from text "1" to numeric 1 ---
$wsPedidos.Cells.Item(10, 4).Value = 1

# --- Pedidos: REQ-010 (row 11) ---
$wsPedidos.Cells.Item(11, 1).Value = "REQ-010"
$wsPedidos.Cells.Item(11, 2).Value = "30/05/2025 08:55"
$wsPedidos.Cells.Item(11, 3).Value = "Renault"
$wsPedidos.Cells.Item(11, 4).Value = 12
$wsPedidos.Cells.Item(11, 5).Value = "R12-LA-A1"
$wsPedidos.Cells.Item(11, 6).Value = "teste"
$wsPedidos.Cells.Item(11, 7).Value = ""
$wsPedidos.Cells.Item(11, 8).Value = "Pendente"
$wsPedidos.Cells.Item(11, 9).Value = ""
$wsPedidos.Cells.Item(11, 10).Value = ""

# --- Pedidos: REQ-011 (row 12) — RACK arrives as text "12" this time ---
$wsPedidos.Cells.Item(12, 1).Value = "REQ-011"
$wsPedidos.Cells.Item(12, 2).Value = "30/05/2025 08:59"
$wsPedidos.Cells.Item(12, 3).Value = "Renault"
$wsPedidos.Cells.Item(12, 4).NumberFormat = "@"
$wsPedidos.Cells.Item(12, 4).Value = "12"
$wsPedidos.Cells.Item(12, 5).Value = "R12-LA-A1"
$wsPedidos.Cells.Item(12, 6).Value = "washington vieira"
$wsPedidos.Cells.Item(12, 7).Value = ""
$wsPedidos.Cells.Item(12, 8).Value = "Pendente"
$wsPedidos.Cells.Item(12, 9).Value = ""
$wsPedidos.Cells.Item(12, 10).Value = ""

# --- Itens: correct D10 ("seccao") from text "1.0" to numeric 1 ---
$wsItens.Cells.Item(10, 4).Value = 1

# --- Itens: REQ-010 (row 11) ---
$wsItens.Cells.Item(11, 1).Value = "REQ-010"
$wsItens.Cells.Item(11, 2).Value = "ACOR2Z-0.35-GY"
$wsItens.Cells.Item(11, 3).Value = "180EX606941"
$wsItens.Cells.Item(11, 4).Value = 0.35
$wsItens.Cells.Item(11, 5).Value = "GY"
$wsItens.Cells.Item(11, 6).Value = 1

# --- Itens: REQ-011 (row 12) — seccao arrives as text "0.35" this time ---
$wsItens.Cells.Item(12, 1).Value = "REQ-011"
$wsItens.Cells.Item(12, 2).Value = "ACOR2Z-0.35-GY"
$wsItens.Cells.Item(12, 3).Value = "180EX606941"
$wsItens.Cells.Item(12, 4).NumberFormat = "@"
$wsItens.Cells.Item(12, 4).Value = "0.35"
$wsItens.Cells.Item(12, 5).Value = "GY"
$wsItens.Cells.Item(12, 6).Value = 1
